# Correccion a Diebold Mariano y revision de Cap1
#
# Applies corrected Diebold-Mariano / HLN test results to the workbook:
#  - Matriz_Resultados: fix a few win/loss/tie indicator cells
#  - P_valores: updated p-values
#  - Estadisticos_HLN_DM: updated HLN test statistics
#  - Resumen_Modelos: updated summary counts / rates that derive from the above

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Matriz_Resultados
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Matriz_Resultados")
$ws.Range("E2").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("D5").Value = 0

# ---------------------------------------------------------------------------
# P_valores
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("P_valores")
$ws.Range("C2").Value = 0.00009250352293754815
$ws.Range("D2").Value = 0.002617164976163044
$ws.Range("E2").Value = 0.01788968829883553

$ws.Range("B3").Value = 0.00009250352293754815
$ws.Range("D3").Value = 0.535310749788648
$ws.Range("E3").Value = 0.141860918160333

$ws.Range("B4").Value = 0.002617164976163044
$ws.Range("C4").Value = 0.535310749788648
$ws.Range("E4").Value = 0.06647042378444534

$ws.Range("B5").Value = 0.01788968829883553
$ws.Range("C5").Value = 0.141860918160333
$ws.Range("D5").Value = 0.06647042378444534

# ---------------------------------------------------------------------------
# Estadisticos_HLN_DM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Estadisticos_HLN_DM")
$ws.Range("C2").Value = -4.76813476679171
$ws.Range("D2").Value = -3.392592122517068
$ws.Range("E2").Value = -2.559133746217347

$ws.Range("B3").Value = 4.76813476679171
$ws.Range("D3").Value = -0.6298023629298993
$ws.Range("E3").Value = 1.52357720040193

$ws.Range("B4").Value = 3.392592122517068
$ws.Range("C4").Value = 0.6298023629298993
$ws.Range("E4").Value = 1.930993469446415

$ws.Range("B5").Value = 2.559133746217347
$ws.Range("C5").Value = -1.52357720040193
$ws.Range("D5").Value = -1.930993469446415

# ---------------------------------------------------------------------------
# Resumen_Modelos
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resumen_Modelos")
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 1
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "66.7%"
$ws.Range("E2").Style = "Normal"

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 3
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.0%"
$ws.Range("E5").Style = "Normal"
